# The Title / Author / Date paragraphs were each typed as one run per
# "word" (<w:r><w:t>Test</w:t></w:r><w:r><w:t> </w:t></w:r> ...). This
# edit collapses each of those paragraphs down to the single run that a
# normal, un-fragmented typed-in-one-go paragraph would have - the
# visible text does not change at all, only how it is split across runs.
#
# Word COM's Range.Find.Execute(..., Replace:=wdReplaceAll) rewrites the
# matched range's contents into a single run carrying ReplaceWith, so
# running it with FindText == ReplaceWith == the paragraph's own full
# text (scoped to that paragraph's Range) performs exactly that merge.

$d = $word.ActiveDocument

$wdReplaceAll    = 2
$wdFindContinue  = 1

function Merge-ParagraphRuns([string]$styleName, [string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs($i)
        if ($para.Style.NameLocal -eq $styleName) {
            $rng = $para.Range
            $rng.Find.Execute(
                $text,            # FindText
                $false,           # MatchCase
                $false,           # MatchWholeWord
                $false,           # MatchWildcards
                $false,           # MatchSoundsLike
                $false,           # MatchAllWordForms
                $true,            # Forward
                $wdFindContinue,  # Wrap
                $false,           # Format
                $text,            # ReplaceWith
                $wdReplaceAll     # Replace
            )
            break
        }
    }
}

Merge-ParagraphRuns "Title"  "Test 011: Unnumbered theorems work"
Merge-ParagraphRuns "Author" "Emma Cliffe, Skills Centre: MASH, University of Bath"
Merge-ParagraphRuns "Date"   "October 2021"
